# Add a new "2022-Q4" sheet (right after "总计") with fund holding data,
# and update the "总计" summary sheet with the new quarter's row.

$wb = $excel.ActiveWorkbook
$totalSheet = $wb.Worksheets.Item(1)

# The workbook uses one "bold / centered / thin-bordered" direct cell format
# (no named style) for header rows and the index column A on every quarter
# sheet. Grab a cell that already carries it so we can clone the formatting
# with Copy/PasteSpecial instead of trying to rebuild it property by
# property.
$styleSrc = $totalSheet.Range("B1")

# ---------------------------------------------------------------------------
# 1. Insert the new "2022-Q4" worksheet right after the first ("总计") sheet.
# ---------------------------------------------------------------------------
$q4Sheet = $wb.Worksheets.Add($null, $totalSheet)
$q4Sheet.Name = "2022-Q4"

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Count; $i++) {
    $q4Sheet.Cells.Item(1, $i + 2).Value = $headers[$i]
}

$q4Data = @(
    @("217024", "招商安盈债券A",            "40.95", "20.32", "2.11", "0.8640", 2),
    @("014887", "招商安福1年定期开放债券",   "17.78", "33.59", "1.82", "0.3236", 5),
    @("010430", "招商安阳债券A",            "16.90", "20.35", "1.44", "0.2434", 5),
    @("016513", "招商安嘉债券",              "16.27", "20.17", "0.71", "0.1155", 9),
    @("010431", "招商安阳债券C",            "0.11",  "20.35", "1.44", "0.0016", 5),
    @("000706", "中邮多策略灵活配置混合",     "0.01",  "89.04", "4.88", "0.0005", 3),
    @("012233", "招商安盈债券C",            "0.01",  "20.32", "2.11", "0.0002", 2)
)

for ($i = 0; $i -lt $q4Data.Count; $i++) {
    $row = $i + 2
    $item = $q4Data[$i]

    $q4Sheet.Cells.Item($row, 1).Value = $i

    # Columns B, D, E, F, G hold text that looks numeric in the source data
    # (fund codes / percentages kept as strings) - force text so the engine
    # does not silently convert them to numbers.
    $textCols = @(2, 4, 5, 6, 7)
    foreach ($col in $textCols) {
        $c = $q4Sheet.Cells.Item($row, $col)
        $c.NumberFormat = "@"
        $c.Value = [string]$item[$col - 2]
    }

    $q4Sheet.Cells.Item($row, 3).Value = [string]$item[1]
    $q4Sheet.Cells.Item($row, 8).Value = $item[6]
}

# Re-apply the shared "bold/center/border" formatting to the header row and
# the index column, then drop the stray quote-prefix / number-format
# styling that NumberFormat="@" / the text assignments above added to the
# data cells (they should stay plain, unstyled cells).
$styleSrc.Copy()
$q4Sheet.Range("B1:H1").PasteSpecial(-4122)
$q4Sheet.Range("A2:A8").PasteSpecial(-4122)
$q4Sheet.Range("B2:H8").ClearFormats()

# ---------------------------------------------------------------------------
# 2. Update the "总计" sheet: insert a new row for 2022-Q4 right under the
#    header, pushing every existing quarter row down by one.
# ---------------------------------------------------------------------------
$totalSheet.Rows.Item(2).Insert()
$totalSheet.Range("B2:D2").ClearFormats()

$totalSheet.Cells.Item(2, 2).Value = "2022-Q4"
$totalSheet.Cells.Item(2, 3).Value = 7
$totalSheet.Cells.Item(2, 4).Value = 1.55

# Renumber the index column (A), which just counts rows from 0, and
# reapply its formatting (the freshly inserted row 2 doesn't have it yet).
for ($r = 2; $r -le 9; $r++) {
    $totalSheet.Cells.Item($r, 1).Value = $r - 2
}
$styleSrc.Copy()
$totalSheet.Range("A2:A9").PasteSpecial(-4122)
